# Insert a new weekly price record as row 59 (Feria Lagunitas de Puerto Montt -
# Pomelo, Start Ruby / Primera, 2021-12-14). Inserting the row shifts the
# existing rows 59-188 down to 60-189, which is exactly the shift visible in
# the rest of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("59").Insert()

$ws.Range("A59").Value = 4
$ws.Range("B59").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value = "Los Lagos"
$ws.Range("D59").Value = 44544
$ws.Range("E59").Value = 10
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100102
$ws.Range("H59").Value = "Cítricos"
$ws.Range("I59").Value = 100102006
$ws.Range("J59").Value = "Pomelo"
$ws.Range("K59").Value = "Start Ruby"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 200
$ws.Range("N59").Value = 11000
$ws.Range("O59").Value = 12000
$ws.Range("P59").Value = 11500
$ws.Range("Q59").Value = "$/caja 14 kilos empedrada"
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 821
$ws.Range("T59").Value = 14
